$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("B2").Value = 50000
$ws.Range("D2").Value = 0.4240605417690221
$ws.Range("E2").Value = 2.826207004934611
$ws.Range("F2").Value = 0.5600000000000001

# Update row 3 values
$ws.Range("B3").Value = 50000
$ws.Range("D3").Value = 0.549008061904701
$ws.Range("E3").Value = 2.739973608910165
$ws.Range("F3").Value = 1.0053

# Remove rows 4 through 10 entirely
$ws.Range("A4:H10").EntireRow.Delete()
